$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 502.23077
$ws.Range("I33").Value = 411.72726
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 411.72726
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -182.72726
$ws.Range("N33").Value = -1458

$ws.Range("H40").Value = 966.5926
$ws.Range("I40").Value = 924.8333
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 924.8333
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -749.8333
$ws.Range("N40").Value = -1350

$ws.Range("H116").Value = 2425
$ws.Range("I116").Value = 1837.5
$ws.Range("J116").Value = 2542.5
$ws.Range("K116").Value = 1837.5
$ws.Range("L116").Value = 2542.5
$ws.Range("M116").Value = 1604.5
$ws.Range("N116").Value = -9426.5

$ws.Range("H137").Value = 2738.2354
$ws.Range("I137").Value = 2441.35
$ws.Range("K137").Value = 7324.049999999999
$ws.Range("M137").Value = -4774.049999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 9999.5
$ws.Range("J43").Value = 9999.5
$ws.Range("L43").Value = 9999.5
$ws.Range("N43").Value = -10625.5

$ws.Range("H86").Value = 16205
$ws.Range("I86").Value = 6000
$ws.Range("J86").Value = 26410
$ws.Range("K86").Value = 6000
$ws.Range("L86").Value = 26410
$ws.Range("M86").Value = -4814
$ws.Range("N86").Value = -28782

$ws.Range("H89").Value = 16205
$ws.Range("I89").Value = 6000
$ws.Range("J89").Value = 26410
$ws.Range("K89").Value = 18000
$ws.Range("L89").Value = 79230
$ws.Range("M89").Value = -12072
$ws.Range("N89").Value = -91086

$ws.Range("H97").Value = 2248.2354
$ws.Range("I97").Value = 946
$ws.Range("J97").Value = 6480.5
$ws.Range("K97").Value = 946
$ws.Range("L97").Value = 6480.5
$ws.Range("M97").Value = -450
$ws.Range("N97").Value = -7472.5

$ws.Range("H102").Value = 6800
$ws.Range("I102").Value = 6800
$ws.Range("K102").Value = 6800
$ws.Range("M102").Value = -5178

$ws.Range("H110").Value = 1040.8
$ws.Range("I110").Value = 922.75
$ws.Range("J110").Value = 1513
$ws.Range("K110").Value = 922.75
$ws.Range("L110").Value = 1513
$ws.Range("M110").Value = 1122.25
$ws.Range("N110").Value = -5603

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2982.5
$ws.Range("I94").Value = 643.3333
$ws.Range("K94").Value = 643.3333
$ws.Range("M94").Value = -192.3333

$ws.Range("H105").Value = 3408.75
$ws.Range("I105").Value = 3741.5386
$ws.Range("J105").Value = 1966.6666
$ws.Range("K105").Value = 3741.5386
$ws.Range("L105").Value = 1966.6666
$ws.Range("M105").Value = -1994.5386
$ws.Range("N105").Value = -5460.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38341.867
$ws.Range("I31").Value = 65564.414
$ws.Range("J31").Value = 2743.1538
$ws.Range("K31").Value = 65564.414
$ws.Range("L31").Value = 2743.1538
$ws.Range("M31").Value = -65269.414
$ws.Range("N31").Value = -3333.1538

$ws.Range("H34").Value = 38341.867
$ws.Range("I34").Value = 65564.414
$ws.Range("J34").Value = 2743.1538
$ws.Range("K34").Value = 65564.414
$ws.Range("L34").Value = 2743.1538
$ws.Range("M34").Value = -65362.414
$ws.Range("N34").Value = -3147.1538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 680.18604
$ws.Range("I5").Value = 448.58334
$ws.Range("J5").Value = 972.7368
$ws.Range("K5").Value = 1345.75002
$ws.Range("L5").Value = 2918.2104
$ws.Range("M5").Value = -1233.75002
$ws.Range("N5").Value = -3142.2104

$ws.Range("H113").Value = 800.6863
$ws.Range("I113").Value = 1067.3462
$ws.Range("J113").Value = 523.36
$ws.Range("K113").Value = 3202.0386
$ws.Range("L113").Value = 1570.08
$ws.Range("M113").Value = -1032.0386
$ws.Range("N113").Value = -5910.08

$ws.Range("H135").Value = 680.18604
$ws.Range("I135").Value = 448.58334
$ws.Range("J135").Value = 972.7368
$ws.Range("K135").Value = 4037.25006
$ws.Range("L135").Value = 8754.6312
$ws.Range("M135").Value = -1502.25006
$ws.Range("N135").Value = -13824.6312

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4884.967
$ws.Range("I70").Value = 4895.1875
$ws.Range("J70").Value = 4873.2856
$ws.Range("K70").Value = 4895.1875
$ws.Range("L70").Value = 4873.2856
$ws.Range("M70").Value = -4625.1875
$ws.Range("N70").Value = -5413.2856

$ws.Range("H73").Value = 4884.967
$ws.Range("I73").Value = 4895.1875
$ws.Range("J73").Value = 4873.2856
$ws.Range("K73").Value = 4895.1875
$ws.Range("L73").Value = 4873.2856
$ws.Range("M73").Value = -3959.1875
$ws.Range("N73").Value = -6745.2856

$ws.Range("H97").Value = 942.2222
$ws.Range("I97").Value = 954.2857
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 954.2857
$ws.Range("L97").Value = 900
$ws.Range("M97").Value = -458.2857
$ws.Range("N97").Value = -1892

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 514.5
$ws.Range("J22").Value = 519.1667
$ws.Range("L22").Value = 519.1667
$ws.Range("N22").Value = -1109.1667

$ws.Range("H27").Value = 514.5
$ws.Range("J27").Value = 519.1667
$ws.Range("L27").Value = 519.1667
$ws.Range("N27").Value = -733.1667

$ws.Range("H68").Value = 2898.7659
$ws.Range("J68").Value = 2956.4443
$ws.Range("L68").Value = 2956.4443
$ws.Range("N68").Value = -4454.4443

$ws.Range("H71").Value = 2898.7659
$ws.Range("J71").Value = 2956.4443
$ws.Range("L71").Value = 14782.2215
$ws.Range("N71").Value = -22270.2215

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 9994.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 9994.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 9994.25
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -10774.25

$ws.Range("H45").Value = 5000
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 5000
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 5000
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -5982

$ws.Range("H74").Value = 8993.143
$ws.Range("J74").Value = 8993.143
$ws.Range("L74").Value = 8993.143
$ws.Range("N74").Value = -10865.143

$ws.Range("H77").Value = 8993.143
$ws.Range("J77").Value = 8993.143
$ws.Range("L77").Value = 26979.429
$ws.Range("N77").Value = -36339.429
